# Automatic update of files.
#
# 1) Column C ("Förändrad" / last-changed date) is bumped by one day
#    (45183 -> 45184, i.e. 2023-09-14 -> 2023-09-15) for every data row
#    (rows 2 through 108).
# 2) The per-case HYPERLINK formulas for the two newest rows (2 and 3)
#    gain a second HYPERLINK argument with the case's display text
#    (", "A xxxxx-2023"") and the last column (Y), which up to now held
#    the formula only as literal/dead text, is turned into a live
#    formula as well.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- 1) Bump the "Förändrad" date for every data row -----------------
$ws.Range("C2:C108").Value = 45184

# --- 2) Refresh the link formulas for row 2 (A 30779-2023 / VASTERVIK)
$ws.Range("S2").Formula = '=HYPERLINK("https://klasma.github.io/LoggingDetectiveFiles/Logging_VASTERVIK/artfynd/A 30779-2023.xlsx, "A 30779-2023"")'
$ws.Range("T2").Formula = '=HYPERLINK("https://klasma.github.io/LoggingDetectiveFiles/Logging_VASTERVIK/kartor/A 30779-2023.png", "A 30779-2023")'
$ws.Range("U2").Formula = '=HYPERLINK("https://klasma.github.io/LoggingDetectiveFiles/Logging_VASTERVIK/knärot/A 30779-2023.png", "A 30779-2023")'
$ws.Range("V2").Formula = '=HYPERLINK("https://klasma.github.io/LoggingDetectiveFiles/Logging_VASTERVIK/klagomål/A 30779-2023.docx", "A 30779-2023")'
$ws.Range("W2").Formula = '=HYPERLINK("https://klasma.github.io/LoggingDetectiveFiles/Logging_VASTERVIK/klagomålsmail/A 30779-2023.docx", "A 30779-2023")'
$ws.Range("X2").Formula = '=HYPERLINK("https://klasma.github.io/LoggingDetectiveFiles/Logging_VASTERVIK/tillsyn/A 30779-2023.docx", "A 30779-2023")'
$ws.Range("Y2").Formula = '=HYPERLINK("https://klasma.github.io/LoggingDetectiveFiles/Logging_VASTERVIK/tillsynsmail/A 30779-2023.docx", "A 30779-2023")'

# --- 3) Refresh the link formulas for row 3 (A 32298-2023 / MONSTERAS)
$ws.Range("S3").Formula = '=HYPERLINK("https://klasma.github.io/LoggingDetectiveFiles/Logging_MONSTERAS/artfynd/A 32298-2023.xlsx, "A 32298-2023"")'
$ws.Range("T3").Formula = '=HYPERLINK("https://klasma.github.io/LoggingDetectiveFiles/Logging_MONSTERAS/kartor/A 32298-2023.png", "A 32298-2023")'
$ws.Range("U3").Formula = '=HYPERLINK("https://klasma.github.io/LoggingDetectiveFiles/Logging_MONSTERAS/knärot/A 32298-2023.png", "A 32298-2023")'
$ws.Range("V3").Formula = '=HYPERLINK("https://klasma.github.io/LoggingDetectiveFiles/Logging_MONSTERAS/klagomål/A 32298-2023.docx", "A 32298-2023")'
$ws.Range("W3").Formula = '=HYPERLINK("https://klasma.github.io/LoggingDetectiveFiles/Logging_MONSTERAS/klagomålsmail/A 32298-2023.docx", "A 32298-2023")'
$ws.Range("X3").Formula = '=HYPERLINK("https://klasma.github.io/LoggingDetectiveFiles/Logging_MONSTERAS/tillsyn/A 32298-2023.docx", "A 32298-2023")'
$ws.Range("Y3").Formula = '=HYPERLINK("https://klasma.github.io/LoggingDetectiveFiles/Logging_MONSTERAS/tillsynsmail/A 32298-2023.docx", "A 32298-2023")'
